# Auto-generated edit script applying the Seraph_Profits.xlsx diff
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row change group #0
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
# Row change group #1
$ws.Range("H95").Value = 11699.4
$ws.Range("J95").Value = 11699.4
$ws.Range("L95").Value = 11699.4
$ws.Range("N95").Value = -17191.4
# Row change group #2
$ws.Range("H137").Value = 2060.625
$ws.Range("I137").Value = 1697.5
$ws.Range("J137").Value = 2423.75
$ws.Range("K137").Value = 5092.5
$ws.Range("L137").Value = 7271.25
$ws.Range("M137").Value = -2542.5
$ws.Range("N137").Value = -12371.25
# Row change group #3
$ws.Range("H138").Value = 4221.143
$ws.Range("J138").Value = 9991.833000000001
$ws.Range("L138").Value = 29975.499
$ws.Range("N138").Value = -40255.499

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row change group #4
$ws.Range("H30").Value = 1990
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()
# Row change group #5
$ws.Range("H45").Value = 2849.75
$ws.Range("J45").Value = 3133.3333
$ws.Range("L45").Value = 3133.3333
$ws.Range("N45").Value = -3887.3333
# Row change group #6
$ws.Range("H102").Value = 15875530
$ws.Range("I102").Value = 55556252
$ws.Range("K102").Value = 55556252
$ws.Range("M102").Value = -55554630
# Row change group #7
$ws.Range("H112").Value = 63499.668
$ws.Range("J112").Value = 63499.668
$ws.Range("L112").Value = 63499.668
$ws.Range("N112").Value = -66453.66800000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row change group #8
$ws.Range("H20").Value = 1417.4375
$ws.Range("I20").Value = 1506.5
$ws.Range("K20").Value = 1506.5
$ws.Range("M20").Value = -1259.5
# Row change group #9
$ws.Range("H64").Value = 1837.6666
$ws.Range("J64").Value = 2006.5
$ws.Range("L64").Value = 2006.5
$ws.Range("N64").Value = -2456.5
# Row change group #10
$ws.Range("H67").Value = 1837.6666
$ws.Range("J67").Value = 2006.5
$ws.Range("L67").Value = 2006.5
$ws.Range("N67").Value = -3566.5
# Row change group #11
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("M87").ClearContents()
$ws.Range("N87").ClearContents()
# Row change group #12
$ws.Range("H88").Value = 12330.5
$ws.Range("J88").Value = 12330.5
$ws.Range("L88").Value = 12330.5
$ws.Range("N88").Value = -13142.5
# Row change group #13
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("M90").ClearContents()
$ws.Range("N90").ClearContents()
# Row change group #14
$ws.Range("H91").Value = 12330.5
$ws.Range("J91").Value = 12330.5
$ws.Range("L91").Value = 12330.5
$ws.Range("N91").Value = -15138.5
# Row change group #15
$ws.Range("H99").Value = 1295.5
$ws.Range("I99").Value = 1295.5
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1295.5
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 202.5
$ws.Range("N99").ClearContents()
# Row change group #16
$ws.Range("H105").Value = 4633982.5
$ws.Range("J105").Value = 4745.143
$ws.Range("L105").Value = 4745.143
$ws.Range("N105").Value = -8239.143

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row change group #17
$ws.Range("H16").Value = 100003420
$ws.Range("I16").Value = 100003420
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 100003420
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -100003133
$ws.Range("N16").ClearContents()
# Row change group #18
$ws.Range("H99").Value = 10268.97
$ws.Range("J99").Value = 14945.6
$ws.Range("L99").Value = 14945.6
$ws.Range("N99").Value = -17941.6
# Row change group #19
$ws.Range("H107").Value = 16667689
$ws.Range("I107").Value = 29412280
$ws.Range("J107").Value = 1684.6154
$ws.Range("K107").Value = 29412280
$ws.Range("L107").Value = 1684.6154
$ws.Range("M107").Value = -29410360
$ws.Range("N107").Value = -5524.6154
# Row change group #20
$ws.Range("H113").Value = 100003420
$ws.Range("I113").Value = 100003420
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 100003420
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -100001250
$ws.Range("N113").ClearContents()
# Row change group #21
$ws.Range("H122").Value = 1248
$ws.Range("I122").Value = 1214.3334
$ws.Range("J122").Value = 1450
$ws.Range("K122").Value = 3643.0002
$ws.Range("L122").Value = 4350
$ws.Range("M122").Value = -1193.0002
$ws.Range("N122").Value = -9250
# Row change group #22
$ws.Range("H126").Value = 10268.97
$ws.Range("J126").Value = 14945.6
$ws.Range("L126").Value = 44836.8
$ws.Range("N126").Value = -49776.8
# Row change group #23
$ws.Range("H134").Value = 3622.5
$ws.Range("I134").Value = 3246.5
$ws.Range("J134").Value = 3998.5
$ws.Range("K134").Value = 9739.5
$ws.Range("L134").Value = 11995.5
$ws.Range("M134").Value = -7204.5
$ws.Range("N134").Value = -17065.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row change group #24
$ws.Range("H2").Value = 301.69232
$ws.Range("I2").Value = 30
$ws.Range("J2").Value = 422.44446
$ws.Range("K2").Value = 180
$ws.Range("L2").Value = 2534.66676
$ws.Range("M2").Value = -67
$ws.Range("N2").Value = -2760.66676
# Row change group #25
$ws.Range("H33").Value = 55613.777
$ws.Range("J33").Value = 166703.33
$ws.Range("L33").Value = 1000219.98
$ws.Range("N33").Value = -1000785.98
# Row change group #26
$ws.Range("H48").Value = 1173.1818
$ws.Range("I48").Value = 490.4
$ws.Range("K48").Value = 1471.2
$ws.Range("M48").Value = -1221.2
# Row change group #27
$ws.Range("H69").Value = 4375
$ws.Range("I69").Value = 2500
$ws.Range("J69").Value = 5000
$ws.Range("K69").Value = 7500
$ws.Range("L69").Value = 15000
$ws.Range("M69").Value = -6689
$ws.Range("N69").Value = -16622
# Row change group #28
$ws.Range("H72").Value = 4375
$ws.Range("I72").Value = 2500
$ws.Range("J72").Value = 5000
$ws.Range("K72").Value = 22500
$ws.Range("L72").Value = 45000
$ws.Range("M72").Value = -18444
$ws.Range("N72").Value = -53112
# Row change group #29
$ws.Range("H80").Value = 3000
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
# Row change group #30
$ws.Range("H83").Value = 3000
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row change group #31
$ws.Range("H5").Value = 3770.2727
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
# Row change group #32
$ws.Range("H19").Value = 250
$ws.Range("I19").Value = 250
$ws.Range("K19").Value = 250
$ws.Range("M19").Value = 38
# Row change group #33
$ws.Range("H80").Value = 3610.125
$ws.Range("J80").Value = 4499.25
$ws.Range("L80").Value = 4499.25
$ws.Range("N80").Value = -6495.25
# Row change group #34
$ws.Range("H83").Value = 3610.125
$ws.Range("J83").Value = 4499.25
$ws.Range("L83").Value = 22496.25
$ws.Range("N83").Value = -32480.25
# Row change group #35
$ws.Range("H107").Value = 3666.6667
$ws.Range("I107").Value = 500
$ws.Range("K107").Value = 500
$ws.Range("M107").Value = 1420
# Row change group #36
$ws.Range("H126").Value = 2750
$ws.Range("I126").Value = 2500
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 7500
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -5030
$ws.Range("N126").Value = -13940

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row change group #37
$ws.Range("H16").Value = 1500
$ws.Range("I16").Value = 1500
$ws.Range("K16").Value = 1500
$ws.Range("M16").Value = -1330
# Row change group #38
$ws.Range("H46").Value = 343.63635
$ws.Range("I46").Value = 300
$ws.Range("J46").Value = 348
$ws.Range("K46").Value = 300
$ws.Range("L46").Value = 348
$ws.Range("M46").Value = -112
$ws.Range("N46").Value = -724
# Row change group #39
$ws.Range("H82").Value = 144627.86
$ws.Range("I82").Value = 2066
$ws.Range("J82").Value = 999999
$ws.Range("K82").Value = 2066
$ws.Range("L82").Value = 999999
$ws.Range("M82").Value = -1705
$ws.Range("N82").Value = -1000721
# Row change group #40
$ws.Range("H85").Value = 144627.86
$ws.Range("I85").Value = 2066
$ws.Range("J85").Value = 999999
$ws.Range("K85").Value = 2066
$ws.Range("L85").Value = 999999
$ws.Range("M85").Value = -818
$ws.Range("N85").Value = -1002495
# Row change group #41
$ws.Range("H110").Value = 50000
$ws.Range("J110").Value = 50000
$ws.Range("L110").Value = 50000
$ws.Range("N110").Value = -58180

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row change group #42
$ws.Range("H2").Value = 1000250
$ws.Range("J2").Value = 500
$ws.Range("L2").Value = 500
$ws.Range("N2").Value = -724
# Row change group #43
$ws.Range("H49").Value = 235995.33
$ws.Range("I49").Value = 249000
$ws.Range("K49").Value = 249000
$ws.Range("M49").Value = -248770
# Row change group #44
$ws.Range("H104").Value = 23750
$ws.Range("J104").Value = 23750
$ws.Range("L104").Value = 23750
$ws.Range("N104").Value = -30738
# Row change group #45
$ws.Range("H107").Value = 425.5
$ws.Range("I107").Value = 425.5
$ws.Range("K107").Value = 1276.5
$ws.Range("M107").Value = 643.5
